# Generate Report for Handback
#
# The handback file f5507395-4caa-43b0-8370-8b7d914df3b2.md has now been
# handed back in sync with en-US. Update the status / datetime / error
# columns across the Overview, zh-cn and de-de sheets to reflect this.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet --------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("L3").Value = "2017-02-21 10:00:55"
$wsZhCn.Range("R3").Value = ""

# --- de-de sheet --------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("L3").Value = "2017-02-21 10:01:19"
$wsDeDe.Range("R3").Value = ""
